# Generate Report for Handoff
#
# Adds two new "handed-off" source files (5caef700-...png and
# bc77e478-...png) to the localization status report, renames the
# existing markdown source file's GUID (9dd8a24a... -> 099c0405...),
# and refreshes the handoff timestamps / target-file hashes on the
# zh-cn and de-de status sheets. Each new source row also gets an
# "IsDependency" row describing the per-locale dependency artifact.

$wb = $excel.ActiveWorkbook

$repoBase    = "https://github.com/OpenLocalizationTest/oltest/blob/2ae8f0fbff5564218b781faacc09d7c667441a65/e2e"
$cfgUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/4eef04e4e35776ca72b70e0c545bb3390db6a24b/.localization-config"
$zhHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6d34ff413676c39192d4ec8cd859b56d8f750934/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eb6a444809321a624717d1e921ad9989f8004afe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$oldGuid = "9dd8a24a-4dd4-4ca8-a56f-b6a3965886a3"
$newGuid = "099c0405-346b-476e-b7fa-2b643c4dc928"
$oldHash = "73519fdf91340b84ead8c1652dd0e2a0f2cdf09c"
$newHash = "3d392ce9d401cc290d5b28ee7968bcb5174283f4"

$mdName      = "$newGuid.md"
$png1Name    = "5caef700-1ca8-4716-994d-20b401be3b18.png"
$png2Name    = "bc77e478-2206-4001-a0bc-848cf0c6f0bc.png"
$cfgName     = ".localization-config"

$zhXlfName   = "$newGuid.$newHash.zh-cn.xlf"
$deXlfName   = "$newGuid.$newHash.de-de.xlf"
$png1DepName = "d59a0e5ca57325b044e3a54d40cbfa255b8f7667.png"
$png2DepName = "8645ff7a0b4d3e8727024eae975c8ef038ad3938.png"
$depDisplay  = "e2e\$mdName"

$readyForHandoff = "Ready for handoff"
$notLocalized    = "Not to be localized"
$include         = "Include"
$isDependency    = "IsDependency"
$ignored         = "Ignored"
$epoch           = "0001-01-01 00:00:00"

$zhTime = "2016-03-10 19:04:08"
$deTime = "2016-03-10 19:04:13"

function Set-HyperlinkCell {
    param($ws, [string]$cellRef, [string]$text, [string]$url)

    $ws.Range($cellRef).Value = $text
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $text) | Out-Null
    $ws.Range($cellRef).Style = "HyperLink"
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-HyperlinkCell $wsOverview "A2" $mdName "$repoBase/$mdName"

Set-HyperlinkCell $wsOverview "A3" $png1Name "$repoBase/$png1Name"
$wsOverview.Range("B3").Value = $readyForHandoff
$wsOverview.Range("C3").Value = $readyForHandoff

Set-HyperlinkCell $wsOverview "A4" $png2Name "$repoBase/$png2Name"
$wsOverview.Range("B4").Value = $readyForHandoff
$wsOverview.Range("C4").Value = $readyForHandoff

Set-HyperlinkCell $wsOverview "A5" $cfgName $cfgUrl
$wsOverview.Range("B5").Value = $notLocalized
$wsOverview.Range("C5").Value = $notLocalized

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-HyperlinkCell $wsZh "A2" $mdName "$repoBase/$mdName"
$wsZh.Range("C2").Value = $zhXlfName
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "$zhHandoffBase/$zhXlfName", "", "", $zhXlfName) | Out-Null
$wsZh.Range("C2").Style = "HyperLink"
$wsZh.Range("D2").Value = $zhTime
$wsZh.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G2").Value = $epoch
$wsZh.Range("H2").Value = $include

Set-HyperlinkCell $wsZh "A3" $png1Name "$repoBase/$png1Name"
$wsZh.Range("B3").Value = $readyForHandoff
Set-HyperlinkCell $wsZh "C3" $png1DepName "$zhHandoffBase/$png1DepName"
$wsZh.Range("D3").Value = $zhTime
$wsZh.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G3").Value = $epoch
$wsZh.Range("H3").Value = $isDependency
$wsZh.Range("I3").Value = $depDisplay

Set-HyperlinkCell $wsZh "A4" $png2Name "$repoBase/$png2Name"
$wsZh.Range("B4").Value = $readyForHandoff
Set-HyperlinkCell $wsZh "C4" $png2DepName "$zhHandoffBase/$png2DepName"
$wsZh.Range("D4").Value = $zhTime
$wsZh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G4").Value = $epoch
$wsZh.Range("H4").Value = $isDependency
$wsZh.Range("I4").Value = $depDisplay

Set-HyperlinkCell $wsZh "A5" $cfgName $cfgUrl
$wsZh.Range("B5").Value = $notLocalized
$wsZh.Range("D5").Value = $epoch
$wsZh.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G5").Value = $epoch
$wsZh.Range("H5").Value = $ignored

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Set-HyperlinkCell $wsDe "A2" $mdName "$repoBase/$mdName"
$wsDe.Range("C2").Value = $deXlfName
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "$deHandoffBase/$deXlfName", "", "", $deXlfName) | Out-Null
$wsDe.Range("C2").Style = "HyperLink"
$wsDe.Range("D2").Value = $deTime
$wsDe.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G2").Value = $epoch
$wsDe.Range("H2").Value = $include

Set-HyperlinkCell $wsDe "A3" $png1Name "$repoBase/$png1Name"
$wsDe.Range("B3").Value = $readyForHandoff
Set-HyperlinkCell $wsDe "C3" $png1DepName "$deHandoffBase/$png1DepName"
$wsDe.Range("D3").Value = $deTime
$wsDe.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G3").Value = $epoch
$wsDe.Range("H3").Value = $isDependency
$wsDe.Range("I3").Value = $depDisplay

Set-HyperlinkCell $wsDe "A4" $png2Name "$repoBase/$png2Name"
$wsDe.Range("B4").Value = $readyForHandoff
Set-HyperlinkCell $wsDe "C4" $png2DepName "$deHandoffBase/$png2DepName"
$wsDe.Range("D4").Value = $deTime
$wsDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G4").Value = $epoch
$wsDe.Range("H4").Value = $isDependency
$wsDe.Range("I4").Value = $depDisplay

Set-HyperlinkCell $wsDe "A5" $cfgName $cfgUrl
$wsDe.Range("B5").Value = $notLocalized
$wsDe.Range("D5").Value = $epoch
$wsDe.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G5").Value = $epoch
$wsDe.Range("H5").Value = $ignored

Write-Output "Generate Report for Handoff: done"
